$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price (column D) and 1h volume change (column E)
# figures with the latest scraped snapshot.

$ws.Range("D2").Value = "26.790.65"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "1.550.27"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.78"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.246"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.37"
$ws.Range("E9").Value = "  -4.12%  "
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "1.768.46"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").Value = "1.554.47"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.68"
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").Value = "26.785.82"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.98"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "213.72"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.26"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("E22").Value = "  -1.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.06"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.70"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.82"
$ws.Range("E27").Value = "  -1.01%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").Value = "1.356.78"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.911"
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.802"
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.991"
$ws.Range("E42").Value = "  -1.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.53"
$ws.Range("E43").Value = "  +3.31%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.76"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.89"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.30"
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("D48").Value = "1.683.16"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.79"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("D51").Value = "0.0₇0965"
$ws.Range("E51").Value = "  -2.46%  "
